# Refresh the standings sheet with the latest stats pulled from the db
# (ranks, totals and changes all shift row-to-row as teams move up/down).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1.0
$ws.Range("B2").Value = 'Lundo’s Legends'
$ws.Range("C2").Value = 139.5
$ws.Range("D2").Value = 62.5
$ws.Range("E2").Value = 77.0
$ws.Range("G2").Value = 15.0
$ws.Range("H2").Value = 2.0
$ws.Range("I2").Value = 13.0

$ws.Range("A3").Value = 2.0
$ws.Range("B3").Value = 'EL Onće'
$ws.Range("C3").Value = 124.5
$ws.Range("D3").Value = 60.5
$ws.Range("E3").Value = 64.0
$ws.Range("G3").Value = 34.0
$ws.Range("H3").Value = 20.0
$ws.Range("I3").Value = 14.0

$ws.Range("A4").Value = 3.0
$ws.Range("B4").Value = 'rainmaker'
$ws.Range("C4").Value = 107.0
$ws.Range("D4").Value = 45.0
$ws.Range("E4").Value = 62.0
$ws.Range("G4").Value = 20.5
$ws.Range("H4").Value = 9.5
$ws.Range("I4").Value = 11.0

$ws.Range("A5").Value = 4.0
$ws.Range("B5").Value = 'Samsquanches'
$ws.Range("C5").Value = 103.5
$ws.Range("D5").Value = 56.0
$ws.Range("E5").Value = 47.5
$ws.Range("G5").Value = 5.5
$ws.Range("H5").Value = 0.0
$ws.Range("I5").Value = 5.5

$ws.Range("A6").Value = 5.0
$ws.Range("B6").Value = 'Epic7'
$ws.Range("C6").Value = 96.0
$ws.Range("D6").Value = 38.0
$ws.Range("E6").Value = 58.0
$ws.Range("G6").Value = 5.0
$ws.Range("H6").Value = -7.0
$ws.Range("I6").Value = 12.0

$ws.Range("A7").Value = 6.5
$ws.Range("B7").Value = 'GOD WILLS IT'
$ws.Range("C7").Value = 89.5
$ws.Range("D7").Value = 47.5
$ws.Range("E7").Value = 42.0
$ws.Range("G7").Value = 5.5
$ws.Range("H7").Value = 0.0
$ws.Range("I7").Value = 5.5

$ws.Range("A8").Value = 6.5
$ws.Range("B8").Value = 'Splitfinger Skadoosh'
$ws.Range("C8").Value = 89.5
$ws.Range("D8").Value = 35.5
$ws.Range("E8").Value = 54.0
$ws.Range("G8").Value = 6.0
$ws.Range("H8").Value = -4.0
$ws.Range("I8").Value = 10.0

$ws.Range("A9").Value = 8.0
$ws.Range("B9").Value = 'confusion'
$ws.Range("C9").Value = 81.5
$ws.Range("D9").Value = 49.5
$ws.Range("E9").Value = 32.0
$ws.Range("G9").Value = 12.5
$ws.Range("H9").Value = 11.5
$ws.Range("I9").Value = 1.0

$ws.Range("A10").Value = 9.0
$ws.Range("B10").Value = 'Swampnuts'
$ws.Range("C10").Value = 79.0
$ws.Range("D10").Value = 35.5
$ws.Range("E10").Value = 43.5
$ws.Range("G10").Value = 5.0
$ws.Range("H10").Value = -3.0
$ws.Range("I10").Value = 8.0

$ws.Range("A11").Value = 10.0
$ws.Range("B11").Value = 'MillerTime'
$ws.Range("C11").Value = 76.0
$ws.Range("D11").Value = 38.5
$ws.Range("E11").Value = 37.5
$ws.Range("G11").Value = -8.0
$ws.Range("H11").Value = -11.0
$ws.Range("I11").Value = 3.0

$ws.Range("A12").Value = 11.0
$ws.Range("B12").Value = 'SmokeWalkers'
$ws.Range("C12").Value = 74.5
$ws.Range("D12").Value = 39.5
$ws.Range("E12").Value = 35.0
$ws.Range("G12").Value = -14.0
$ws.Range("H12").Value = -23.0
$ws.Range("I12").Value = 9.0

$ws.Range("A13").Value = 12.0
$ws.Range("B13").Value = 'Mac'
$ws.Range("C13").Value = 74.0
$ws.Range("D13").Value = 38.0
$ws.Range("E13").Value = 36.0
$ws.Range("G13").Value = 4.5
$ws.Range("H13").Value = 2.5
$ws.Range("I13").Value = 2.0

$ws.Range("A14").Value = 13.0
$ws.Range("B14").Value = 'DJ''s Quality Team'
$ws.Range("C14").Value = 64.5
$ws.Range("D14").Value = 40.5
$ws.Range("E14").Value = 24.0
$ws.Range("G14").Value = 8.0
$ws.Range("H14").Value = 2.5
$ws.Range("I14").Value = 5.5

$ws.Range("A15").Value = 14.0
$ws.Range("B15").Value = 'Corbin Copy'
$ws.Range("C15").Value = 61.0
$ws.Range("D15").Value = 43.5
$ws.Range("E15").Value = 17.5
$ws.Range("G15").Value = 5.5
$ws.Range("H15").Value = 0.0
$ws.Range("I15").Value = 5.5
